$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "29.524.35"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.923.79"
$ws.Range("E3").Value = "  +0.70%  "
Set-TextValue "D4" "1.013"
$ws.Range("E4").Value = "  +0.52%  "
Set-TextValue "D5" "326.23"
$ws.Range("E5").Value = "  +0.22%  "
Set-TextValue "D6" "1.012"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  -0.92%  "
Set-TextValue "D8" "0.4045"
$ws.Range("E8").Value = "  -0.52%  "
Set-TextValue "D9" "0.08183"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  -0.57%  "
Set-TextValue "D11" "23.82"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").Value = "1.947.53"
$ws.Range("E12").Value = "  +0.53%  "
Set-TextValue "D13" "6.087"
$ws.Range("E13").Value = "  +1.10%  "
Set-TextValue "D14" "7.311"
$ws.Range("E14").Value = "  +1.89%  "
Set-TextValue "D15" "91.53"
$ws.Range("E15").Value = "  +1.19%  "
Set-TextValue "D16" "0.06889"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("E17").Value = "  +0.53%  "
Set-TextValue "D18" "0.00001038"
$ws.Range("E18").Value = "  +0.10%  "
Set-TextValue "D20" "1.010"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "29.512.35"
$ws.Range("E21").Value = "  -0.01%  "
Set-TextValue "D22" "5.658"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("E23").Value = "  +1.70%  "
Set-TextValue "D24" "2.175"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").Value = "2.148.17"
$ws.Range("E25").Value = "  -0.90%  "
Set-TextValue "D26" "155.91"
$ws.Range("E26").Value = "  -0.75%  "
Set-TextValue "D27" "6.387"
$ws.Range("E27").Value = "  -2.51%  "
Set-TextValue "D28" "20.01"
$ws.Range("E28").Value = "  -0.46%  "
Set-TextValue "D29" "2.083"
$ws.Range("E29").Value = "  -1.63%  "
Set-TextValue "D30" "120.53"
$ws.Range("E30").Value = "  +0.03%  "
Set-TextValue "D31" "1.012"
$ws.Range("E31").Value = "  -1.71%  "
Set-TextValue "D32" "0.09585"
$ws.Range("E32").Value = "  +0.51%  "
Set-TextValue "D33" "5.594"
$ws.Range("E33").Value = "  +1.38%  "
Set-TextValue "D34" "3.563"
$ws.Range("E34").Value = "  +0.07%  "
Set-TextValue "D35" "1.384"
$ws.Range("E35").Value = "  -0.62%  "
Set-TextValue "D36" "0.06353"
$ws.Range("E36").Value = "  +3.80%  "
Set-TextValue "D37" "0.02278"
$ws.Range("E37").Value = "  +0.01%  "
Set-TextValue "D38" "1.190"
Set-TextValue "D39" "0.5934"
$ws.Range("E39").Value = "  -0.57%  "
Set-TextValue "D40" "10.73"
$ws.Range("E40").Value = "  -1.32%  "
Set-TextValue "D41" "1.011"
$ws.Range("E41").Value = "  +0.48%  "
Set-TextValue "D42" "7.888"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("E43").Value = "  -0.75%  "
Set-TextValue "D44" "2.476"
$ws.Range("E44").Value = "  +4.00%  "
Set-TextValue "D45" "1.243"
$ws.Range("E45").Value = "  -2.76%  "
Set-TextValue "D46" "12.33"
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("E47").Value = "  -1.98%  "
Set-TextValue "D48" "0.5540"
$ws.Range("E48").Value = "  -0.73%  "
Set-TextValue "D49" "1.969"
$ws.Range("E49").Value = "  +0.87%  "
Set-TextValue "D50" "117.72"
$ws.Range("E50").Value = "  +0.87%  "
Set-TextValue "D51" "2.433"
$ws.Range("E51").Value = "  +1.06%  "
